$d = $word.ActiveDocument

for ($i = 1; $i -le 13; $i++) {
    try {
        $story = $d.StoryRanges.Item($i)
        if ($story -ne $null) {
            Write-Output ("Story idx " + $i + " Type=" + $story.StoryType + " InlineShapes.Count=" + $story.InlineShapes.Count)
        } else {
            Write-Output ("Story idx " + $i + " is null")
        }
    } catch {
        Write-Output ("Story idx " + $i + " ERROR: " + $_.Exception.Message)
    }
}
